$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The source site re-scraped / re-sorted several same-day fixtures, so the
# match rows for those dates need their F:V (home..url_partida) content
# rotated among themselves. Columns A:E (Indice, pais, torneio, temporada,
# data_partida) stay untouched.
#
# Mapping: destination row -> source row (value to copy FROM, read BEFORE
# any writes happen so overlapping/cyclic swaps don't clobber data).
# ---------------------------------------------------------------------------
$rowMap = @{
    2  = 3
    3  = 2
    5  = 7
    6  = 5
    7  = 6
    13 = 14
    14 = 13
    26 = 27
    27 = 26
    33 = 34
    34 = 35
    35 = 33
    43 = 45
    44 = 43
    45 = 44
    51 = 52
    52 = 51
    55 = 57
    56 = 55
    57 = 56
    79 = 80
    80 = 79
}

# Snapshot F:V for every row referenced above (as source) before writing.
$snapshot = @{}
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("F" + $srcRow + ":V" + $srcRow).Value()
    }
}

# Now write the snapshotted source data into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("F" + $destRow + ":V" + $destRow).Value = $snapshot[$srcRow]
}

# ---------------------------------------------------------------------------
# Append 3 new match rows (106, 107, 108) at the bottom of the table.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=106; A=105; E=45263.625;         F="Koge";       G=1; H="Aalborg";   I=2;
       J=6.03; K="26/11/2023 15:12"; L=9.69; M="03/12/2023 14:52";
       N=4.89; O="26/11/2023 15:12"; P=5.94; Q="03/12/2023 14:52";
       R=1.45; S="26/11/2023 15:12"; T=1.27; U="03/12/2023 14:35";
       V="https://www.betexplorer.com/football/denmark/1st-division/koge-aalborg/jJWUBLBA/" },
    @{ Row=107; A=106; E=45265.77083333334; F="Kolding IF"; G=0; H="B.93";      I=1;
       J=1.38; K="25/11/2023 18:13"; L=1.55; M="05/12/2023 18:28";
       N=5.02; O="25/11/2023 18:13"; P=4.26; Q="05/12/2023 18:28";
       R=7.22; S="25/11/2023 18:13"; T=5.88; U="05/12/2023 18:28";
       V="https://www.betexplorer.com/football/denmark/1st-division/kolding-if-boldklubben-1893/U9VYAuRG/" },
    @{ Row=108; A=107; E=45265.77083333334; F="Naestved";   G=2; H="Helsingor"; I=2;
       J=1.9;  K="27/11/2023 07:42"; L=2.07; M="05/12/2023 18:16";
       N=3.86; O="27/11/2023 07:42"; P=3.4;  Q="05/12/2023 18:29";
       R=3.46; S="27/11/2023 07:42"; T=3.73; U="05/12/2023 18:16";
       V="https://www.betexplorer.com/football/denmark/1st-division/naestved-if-helsingor/0UwIENRi/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A" + $r).Value = $nr.A
    $ws.Range("B" + $r).Value = "denmark"
    $ws.Range("C" + $r).Value = "1st-division"
    $ws.Range("D" + $r).Value = "2023-2024"
    $ws.Range("E" + $r).Value = $nr.E
    $ws.Range("F" + $r).Value = $nr.F
    $ws.Range("G" + $r).Value = $nr.G
    $ws.Range("H" + $r).Value = $nr.H
    $ws.Range("I" + $r).Value = $nr.I
    $ws.Range("J" + $r).Value = $nr.J
    $ws.Range("K" + $r).Value = $nr.K
    $ws.Range("L" + $r).Value = $nr.L
    $ws.Range("M" + $r).Value = $nr.M
    $ws.Range("N" + $r).Value = $nr.N
    $ws.Range("O" + $r).Value = $nr.O
    $ws.Range("P" + $r).Value = $nr.P
    $ws.Range("Q" + $r).Value = $nr.Q
    $ws.Range("R" + $r).Value = $nr.R
    $ws.Range("S" + $r).Value = $nr.S
    $ws.Range("T" + $r).Value = $nr.T
    $ws.Range("U" + $r).Value = $nr.U
    $ws.Range("V" + $r).Value = $nr.V

    # Replicate formatting of column A (bold/border/center) and column E
    # (date number format) from the last existing data row (105).
    $ws.Range("A105").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Range("E105").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
